$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to match the workbook title
$ws.Name = "Gamecube controller data layout"

# The "Powered" label was documented in the wrong column (F3); move it to G3
$ws.Range("G3").Value = "Powered"
$ws.Range("F3").ClearContents()

# Update the active cell/selection saved with the sheet view
$ws.Activate()
$ws.Range("G3").Select()
